# edit.ps1 - apply the protocol_definition_hr.docx changes:
#  1. ENCRYPTION PUB KEY paragraph: collapse the "2048"/"270" split runs
#     back into the surrounding sentence (pure text/run cleanup, wording
#     unchanged).
#  2. SENDER SIGNING KEY paragraph: "nam je dostavljen" -> "smo dostavili"
#     (we delivered it, instead of: it was delivered to us).
#  3. Normal style: turn off "allow punctuation to extend past the text
#     margin" (w:overflowPunct true -> false), exposed on the Word object
#     model as ParagraphFormat.HangingPunctuation.

$d = $word.ActiveDocument

# --- 1. ENCRYPTION PUB KEY: merge the fragmented runs back into one
#        continuous sentence. Scope the Find to the owning paragraph so we
#        can't accidentally touch the same words elsewhere (e.g. the
#        summary table).
$encKeyOld = "ENCRYPTION PUB KEY javni je ključ tipa RSA 2048 bita. " + `
    "Također ga pohranjuje primatelj, koristi ga kako bi enkriptirao " + `
    "sadržaj poruke koji šalje. Odnosno kako bi enkriptirao simetrični " + `
    "ključ koji se koristio za enkriptiranje poruke. Ključ je enkodiran " + `
    "u DER obliku, i kao takav javni ključ je dug 270 bajta (testirano)."

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "ENCRYPTION PUB KEY javni je*") {
        $r = $para.Range
        $found = $r.Find.Execute($encKeyOld, $true, $false, $false, $false, `
            $false, $true, 1, $false, $encKeyOld, 2)
        break
    }
}
Write-Host ("encryption pub key paragraph updated: " + $found)

# --- 2. SENDER SIGNING KEY: reword "koji nam je dostavljen" -> "koji smo
#        dostavili". Again scope to the specific paragraph (there are two
#        "SENDER SIGNING KEY ..." paragraphs in the doc).
$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "SENDER SIGNING KEY javni*nam je dostavljen*") {
        $r = $para.Range
        $found2 = $r.Find.Execute("nam je dostavljen", $true, $false, $false, `
            $false, $false, $true, 1, $false, "smo dostavili", 2)
        break
    }
}
Write-Host ("sender signing key paragraph updated: " + $found2)

# --- 3. Normal style: w:overflowPunct val="true" -> val="false"
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $false
Write-Host ("HangingPunctuation now: " + $normal.ParagraphFormat.HangingPunctuation)
